# Auto update Excel log
# Appends new sensor-log rows (2026-02-06, ~10:10-10:11) to the PIR,
# Humidity and Temperature sheets, matching the upstream data feed.

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $r, $dateVal, $timeVal, $hourVal, $locVal, $valueVal, $statusVal) {
    # Column A holds a plain "YYYY-MM-DD" text date and column E sometimes
    # holds a plain "xx.x%" reading. Excel's automatic type-detection would
    # otherwise convert those into a date serial number / percentage number,
    # so force Text formatting before assigning the value, then restore the
    # default "Normal" style so the cell keeps looking unformatted, matching
    # the rest of the log.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateVal
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $timeVal
    $ws.Cells.Item($r, 3).Value = $hourVal
    $ws.Cells.Item($r, 4).Value = $locVal

    $valueCell = $ws.Cells.Item($r, 5)
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $valueVal
    $valueCell.Style = "Normal"

    $ws.Cells.Item($r, 6).Value = $statusVal
}

# --- PIR sheet: rows 366-378 (Bathroom / No Motion / Inactive) ---
$pirSheet = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(366, "2026-02-06", "10:10:25", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(367, "2026-02-06", "10:10:26", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(368, "2026-02-06", "10:10:29", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(369, "2026-02-06", "10:10:34", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(370, "2026-02-06", "10:10:39", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(371, "2026-02-06", "10:10:44", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(372, "2026-02-06", "10:10:49", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(373, "2026-02-06", "10:10:54", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(374, "2026-02-06", "10:10:59", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(375, "2026-02-06", "10:11:04", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(376, "2026-02-06", "10:11:09", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(377, "2026-02-06", "10:11:14", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(378, "2026-02-06", "10:11:19", "10:00", "Bathroom", "No Motion", "Inactive")
)
foreach ($row in $pirRows) {
    Add-LogRow $pirSheet $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $row[6]
}

# --- Humidity sheet: rows 247-259 (Bathroom / xx.x% / Active) ---
$humiditySheet = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(247, "2026-02-06", "10:10:24", "10:00", "Bathroom", "68.8%", "Active"),
    @(248, "2026-02-06", "10:10:27", "10:00", "Bathroom", "69.0%", "Active"),
    @(249, "2026-02-06", "10:10:32", "10:00", "Bathroom", "68.8%", "Active"),
    @(250, "2026-02-06", "10:10:37", "10:00", "Bathroom", "68.9%", "Active"),
    @(251, "2026-02-06", "10:10:42", "10:00", "Bathroom", "68.0%", "Active"),
    @(252, "2026-02-06", "10:10:47", "10:00", "Bathroom", "69.0%", "Active"),
    @(253, "2026-02-06", "10:10:53", "10:00", "Bathroom", "68.1%", "Active"),
    @(254, "2026-02-06", "10:10:57", "10:00", "Bathroom", "68.9%", "Active"),
    @(255, "2026-02-06", "10:11:03", "10:00", "Bathroom", "68.0%", "Active"),
    @(256, "2026-02-06", "10:11:07", "10:00", "Bathroom", "68.8%", "Active"),
    @(257, "2026-02-06", "10:11:13", "10:00", "Bathroom", "68.9%", "Active"),
    @(258, "2026-02-06", "10:11:17", "10:00", "Bathroom", "68.8%", "Active"),
    @(259, "2026-02-06", "10:11:23", "10:00", "Bathroom", "68.8%", "Active")
)
foreach ($row in $humidityRows) {
    Add-LogRow $humiditySheet $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $row[6]
}

# --- Temperature sheet: rows 247-259 (Bathroom / xx.xC / Active) ---
$temperatureSheet = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(247, "2026-02-06", "10:10:25", "10:00", "Bathroom", "28.0C", "Active"),
    @(248, "2026-02-06", "10:10:28", "10:00", "Bathroom", "28.1C", "Active"),
    @(249, "2026-02-06", "10:10:33", "10:00", "Bathroom", "28.0C", "Active"),
    @(250, "2026-02-06", "10:10:38", "10:00", "Bathroom", "28.1C", "Active"),
    @(251, "2026-02-06", "10:10:43", "10:00", "Bathroom", "28.1C", "Active"),
    @(252, "2026-02-06", "10:10:48", "10:00", "Bathroom", "28.1C", "Active"),
    @(253, "2026-02-06", "10:10:53", "10:00", "Bathroom", "28.1C", "Active"),
    @(254, "2026-02-06", "10:10:58", "10:00", "Bathroom", "28.0C", "Active"),
    @(255, "2026-02-06", "10:11:03", "10:00", "Bathroom", "28.1C", "Active"),
    @(256, "2026-02-06", "10:11:08", "10:00", "Bathroom", "28.0C", "Active"),
    @(257, "2026-02-06", "10:11:13", "10:00", "Bathroom", "28.1C", "Active"),
    @(258, "2026-02-06", "10:11:18", "10:00", "Bathroom", "28.0C", "Active"),
    @(259, "2026-02-06", "10:11:23", "10:00", "Bathroom", "28.0C", "Active")
)
foreach ($row in $temperatureRows) {
    Add-LogRow $temperatureSheet $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $row[6]
}
